# Apply the changes described in the diff.

$wb = $excel.ActiveWorkbook

# --- Astronauta (sheet 1) ---
$wsAstronauta = $wb.Worksheets.Item("Astronauta")
$wsAstronauta.Range("E3").Value = 1
$wsAstronauta.Range("E11").Value = 1
$wsAstronauta.Range("E19").Value = 0
$wsAstronauta.Range("E11").Select()

# --- Mago (sheet 3) ---
$wsMago = $wb.Worksheets.Item("Mago")
$wsMago.Range("E5").Value = 1
$wsMago.Range("E8").Value = 0
$wsMago.Range("E16").Value = 0
$wsMago.Range("E22").Value = 1
$wsMago.Range("E23").Value = 1
$wsMago.Range("E28").Value = 0
$wsMago.Range("A24").Select()

# --- Ninja (sheet 4) ---
$wsNinja = $wb.Worksheets.Item("Ninja")
$wsNinja.Range("E24").Value = 1
$wsNinja.Range("F24").Value = 1
$wsNinja.Activate()
$wsNinja.Range("G25").Select()
$excel.ActiveWindow.Zoom = 100
